$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-26 (columns A-I).
# Column A = class label (string), columns B-I = numeric thresholds.
$data = @(
    @(2,  "L107",  140, 6,   5.33, 2.35, 1.04, 0.64, 0.28, 0),
    @(3,  "L109",  140, 6,   5.33, 2.35, 1.04, 0.64, 0.28, 0),
    @(4,  "L104",  140, 4,   3.79, 1.6,  0.64, 0.4,  0.18, 0),
    @(5,  "L105a", 140, 4,   3.79, 1.6,  0.64, 0.4,  0.18, 0),
    @(6,  "L207",  140, 4,   3.79, 1.6,  0.64, 0.4,  0.18, 0),
    @(7,  "L105b", 140, 3.6, 1.9,  0.77, 0.4,  0.18, 0.11, 0),
    @(8,  "L106",  140, 6,   4.6,  2,    1,    0.6,  0.3,  0),
    @(9,  "L208",  140, 6,   4.6,  2,    1,    0.6,  0.3,  0),
    @(10, "L101",  140, 3,   1.9,  0.77, 0.4,  0.18, 0.11, 0),
    @(11, "L102",  140, 3,   1.9,  0.77, 0.4,  0.18, 0.11, 0),
    @(12, "L201",  140, 3,   1.9,  0.77, 0.4,  0.18, 0.11, 0),
    @(13, "L202",  140, 3,   1.9,  0.77, 0.4,  0.18, 0.11, 0),
    @(14, "L204",  140, 3,   1.9,  0.77, 0.4,  0.18, 0.11, 0),
    @(15, "L205",  140, 3,   1.9,  0.77, 0.4,  0.18, 0.11, 0),
    @(16, "L103",  140, 3.6, 3.46, 1.46, 0.64, 0.4,  0.18, 0),
    @(17, "L203",  140, 3.6, 3.46, 1.46, 0.64, 0.4,  0.18, 0),
    @(18, "L206",  140, 3.6, 3.46, 1.46, 0.64, 0.4,  0.18, 0),
    @(19, "L108",  140, 7,   6.03, 2.66, 1.24, 0.77, 0.34, 0),
    @(20, "L110",  140, 7,   6.03, 2.66, 1.24, 0.77, 0.34, 0),
    @(21, "L301",  140, 3,   1.46, 0.64, 0.23, 0.13, 0.06, 0),
    @(22, "L302",  140, 3,   1.46, 0.64, 0.23, 0.13, 0.06, 0),
    @(23, "L304",  140, 3,   1.46, 0.64, 0.23, 0.13, 0.06, 0),
    @(24, "L305",  140, 3,   1.46, 0.64, 0.23, 0.13, 0.06, 0),
    @(25, "L303",  140, 3,   1,    0.77, 0.4,  0.18, 0.11, 0),
    @(26, "L306",  140, 3,   1,    0.77, 0.4,  0.18, 0.11, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    for ($c = 2; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c]
    }
}
